# Trade #73 closed at 2026-02-18 00:27:37 - unknown UNKNOWN +0.000%
#
# Updates live trading results after trade #101 (row 102 on "All Trades",
# row 34 on "MarketMaking") closed via early_exit, and a new trade #130
# (row 131 / row 51) opened.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.22   # Current Capital
$summary.Range("B4").Value = 0.33      # Total P&L $
$summary.Range("B6").Value = 101       # Total Trades
$summary.Range("B8").Value = 38        # Losing Trades
$summary.Range("B9").Value = 46.53     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.37      # Capital
$status.Range("D6").Value = 33         # Trades
$status.Range("E6").Value = -0.44      # P&L $
$status.Range("F6").Value = -0.63      # P&L %
$status.Range("G6").Value = 45.45      # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Trade #101 (row 102) closes via early_exit
$allTrades.Range("G102").Value = 0.64
$allTrades.Range("H102").Value = "CLOSED"
$allTrades.Range("I102").Value = -4.4776
$allTrades.Range("J102").Value = -0.03
$allTrades.Range("K102").Value = 99.37
$allTrades.Range("L102").Value = "early_exit"
$allTrades.Range("M102").Value = 0.14

# New trade #130 (row 131) opens
$allTrades.Range("A131").Value = 130
$allTrades.Range("B131").NumberFormat = "@"
$allTrades.Range("B131").Value = "2026-02-18"
$allTrades.Range("B131").Style = "Normal"
$allTrades.Range("C131").Value = "00:27:31"
$allTrades.Range("D131").Value = "MarketMaking"
$allTrades.Range("E131").Value = "DOWN"
$allTrades.Range("F131").Value = 0.67
$allTrades.Range("H131").Value = "OPEN"
$allTrades.Range("I131").Value = 0
$allTrades.Range("J131").Value = 0
$allTrades.Range("K131").Value = 99.39967800952272
$allTrades.Range("M131").Value = 0
$allTrades.Range("N131").Value = 0
$allTrades.Range("O131").Value = 0
$allTrades.Range("P131").Value = 0.65
$allTrades.Range("Q131").Value = "Wide spread capture: 392 bps vs avg 291 bps"

# ---------------------------------------------------------------------
# MarketMaking sheet (strategy-specific trade log)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Trade #101 (row 34) closes via early_exit
$mm.Range("G34").Value = 0.64
$mm.Range("H34").Value = "CLOSED"
$mm.Range("I34").Value = -4.4776
$mm.Range("J34").Value = -0.03
$mm.Range("K34").Value = 99.37
$mm.Range("P34").Value = "early_exit"
$mm.Range("Q34").Value = 0.14

# New trade #130 (row 51) opens
$mm.Range("A51").Value = 130
$mm.Range("B51").NumberFormat = "@"
$mm.Range("B51").Value = "2026-02-18"
$mm.Range("B51").Style = "Normal"
$mm.Range("C51").Value = "00:27:31"
$mm.Range("D51").Value = "MarketMaking"
$mm.Range("E51").Value = "DOWN"
$mm.Range("F51").Value = 0.67
$mm.Range("H51").Value = "OPEN"
$mm.Range("I51").Value = 0
$mm.Range("J51").Value = 0
$mm.Range("K51").Value = 99.39967800952272
$mm.Range("L51").Value = 0
$mm.Range("M51").Value = 0
$mm.Range("N51").Value = 0.65
$mm.Range("O51").Value = "Wide spread capture: 392 bps vs avg 291 bps"
$mm.Range("Q51").Value = 0
